# Update collection counts (column F) by +1 for a set of matching events
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row numbers -> target F-column values on "展览"
$wsExhibit.Range("F11").Value = 1381
$wsExhibit.Range("F12").Value = 3030
$wsExhibit.Range("F13").Value = 497
$wsExhibit.Range("F24").Value = 26
$wsExhibit.Range("F25").Value = 3570

# Same events, mirrored on "全部类型"
$wsAll.Range("F22").Value = 1381
$wsAll.Range("F23").Value = 3030
$wsAll.Range("F24").Value = 497
$wsAll.Range("F37").Value = 26
$wsAll.Range("F38").Value = 3570
